$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 254
$ws.Range("J5").Value = 350
$ws.Range("L5").Value = 350
$ws.Range("N5").Value = -580

$ws.Range("H17").Value = 1877
$ws.Range("J17").Value = 1877
$ws.Range("L17").Value = 5631
$ws.Range("N17").Value = -5967

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()

$ws.Range("H62").Value = 10240.15
$ws.Range("I62").Value = 9754.077
$ws.Range("K62").Value = 9754.077
$ws.Range("M62").Value = -9130.077

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("N64").ClearContents()

$ws.Range("H65").Value = 10240.15
$ws.Range("I65").Value = 9754.077
$ws.Range("K65").Value = 48770.38499999999
$ws.Range("M65").Value = -45650.38499999999

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("N67").ClearContents()

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()

$ws.Range("H86").Value = 4469.8237
$ws.Range("I86").Value = 1746.75
$ws.Range("K86").Value = 1746.75
$ws.Range("M86").Value = -623.75

$ws.Range("H89").Value = 4469.8237
$ws.Range("I89").Value = 1746.75
$ws.Range("K89").Value = 8733.75
$ws.Range("M89").Value = -3117.75

$ws.Range("H98").Value = 1711.4762
$ws.Range("I98").Value = 1746.5333
$ws.Range("J98").Value = 1623.8334
$ws.Range("K98").Value = 1746.5333
$ws.Range("L98").Value = 1623.8334
$ws.Range("M98").Value = -248.5333000000001
$ws.Range("N98").Value = -4619.8334

$ws.Range("H106").Value = 2149.1667
$ws.Range("I106").Value = 1543.3334
$ws.Range("K106").Value = 1543.3334
$ws.Range("M106").Value = -912.3334

$ws.Range("H112").Value = 1400.9062
$ws.Range("J112").Value = 1458.4138
$ws.Range("L112").Value = 4375.2414
$ws.Range("N112").Value = -6591.2414

$ws.Range("H116").Value = 38755.637
$ws.Range("I116").Value = 62663.418
$ws.Range("J116").Value = 10066.3
$ws.Range("K116").Value = 62663.418
$ws.Range("L116").Value = 10066.3
$ws.Range("M116").Value = -59221.418
$ws.Range("N116").Value = -16950.3

$ws.Range("H122").Value = 1711.4762
$ws.Range("I122").Value = 1746.5333
$ws.Range("J122").Value = 1623.8334
$ws.Range("K122").Value = 5239.5999
$ws.Range("L122").Value = 4871.5002
$ws.Range("M122").Value = -2789.5999
$ws.Range("N122").Value = -9771.5002

$ws.Range("H132").Value = 1922.079
$ws.Range("I132").Value = 1865.9459
$ws.Range("K132").Value = 5597.8377
$ws.Range("M132").Value = -3067.8377

$ws.Range("H135").Value = 1180.8462
$ws.Range("I135").Value = 1154.4166
$ws.Range("K135").Value = 10389.7494
$ws.Range("M135").Value = -7854.749400000001

$ws.Range("H137").Value = 2642.9565
$ws.Range("I137").Value = 1395.4286
$ws.Range("J137").Value = 3188.75
$ws.Range("K137").Value = 4186.2858
$ws.Range("L137").Value = 9566.25
$ws.Range("M137").Value = -1636.2858
$ws.Range("N137").Value = -14666.25

$ws.Range("H138").Value = 2091.3044
$ws.Range("I138").Value = 1477.6451
$ws.Range("K138").Value = 4432.9353
$ws.Range("M138").Value = 707.0646999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2998
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 2998
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 2998
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -3572

$ws.Range("H45").Value = 2826.8572
$ws.Range("I45").Value = 3519.5881
$ws.Range("J45").Value = 1756.2727
$ws.Range("K45").Value = 3519.5881
$ws.Range("L45").Value = 1756.2727
$ws.Range("M45").Value = -3142.5881
$ws.Range("N45").Value = -2510.2727

$ws.Range("H61").Value = 3847.9678
$ws.Range("I61").Value = 2512
$ws.Range("K61").Value = 2512
$ws.Range("M61").Value = -2300

$ws.Range("H74").Value = 2694.6667
$ws.Range("I74").Value = 1731
$ws.Range("K74").Value = 1731
$ws.Range("M74").Value = -857

$ws.Range("H77").Value = 2694.6667
$ws.Range("I77").Value = 1731
$ws.Range("K77").Value = 8655
$ws.Range("M77").Value = -4287

$ws.Range("H88").Value = 1658.909
$ws.Range("I88").Value = 1343.5
$ws.Range("J88").Value = 2500
$ws.Range("K88").Value = 1343.5
$ws.Range("L88").Value = 2500
$ws.Range("M88").Value = -937.5
$ws.Range("N88").Value = -3312

$ws.Range("H91").Value = 1658.909
$ws.Range("I91").Value = 1343.5
$ws.Range("J91").Value = 2500
$ws.Range("K91").Value = 1343.5
$ws.Range("L91").Value = 2500
$ws.Range("M91").Value = 60.5
$ws.Range("N91").Value = -5308

$ws.Range("H96").Value = 42171.5
$ws.Range("J96").Value = 42171.5
$ws.Range("L96").Value = 42171.5
$ws.Range("N96").Value = -47663.5

$ws.Range("H122").Value = 3111.756
$ws.Range("I122").Value = 2728.6667
$ws.Range("J122").Value = 3850.5715
$ws.Range("K122").Value = 8186.000100000001
$ws.Range("L122").Value = 11551.7145
$ws.Range("M122").Value = -5736.000100000001
$ws.Range("N122").Value = -16451.7145

$ws.Range("H132").Value = 3299.16
$ws.Range("I132").Value = 3299.16
$ws.Range("K132").Value = 9897.48
$ws.Range("M132").Value = -7367.48

$ws.Range("H136").Value = 3847.9678
$ws.Range("I136").Value = 2512
$ws.Range("K136").Value = 7536
$ws.Range("M136").Value = -4986

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1981.5238
$ws.Range("I20").Value = 2155.3
$ws.Range("K20").Value = 2155.3
$ws.Range("M20").Value = -1908.3

$ws.Range("H86").Value = 1458.9
$ws.Range("I86").Value = 1526.1428
$ws.Range("J86").Value = 1302
$ws.Range("K86").Value = 1526.1428
$ws.Range("L86").Value = 1302
$ws.Range("M86").Value = -403.1428000000001
$ws.Range("N86").Value = -3548

$ws.Range("H88").Value = 49959.145
$ws.Range("J88").Value = 49959.145
$ws.Range("L88").Value = 49959.145
$ws.Range("N88").Value = -50771.145

$ws.Range("H89").Value = 1458.9
$ws.Range("I89").Value = 1526.1428
$ws.Range("J89").Value = 1302
$ws.Range("K89").Value = 7630.714
$ws.Range("L89").Value = 6510
$ws.Range("M89").Value = -2014.714
$ws.Range("N89").Value = -17742

$ws.Range("H91").Value = 49959.145
$ws.Range("J91").Value = 49959.145
$ws.Range("L91").Value = 49959.145
$ws.Range("N91").Value = -52767.145

$ws.Range("H105").Value = 1343.6
$ws.Range("I105").Value = 1369.4783
$ws.Range("J105").Value = 1258.5714
$ws.Range("K105").Value = 1369.4783
$ws.Range("L105").Value = 1258.5714
$ws.Range("M105").Value = 377.5217
$ws.Range("N105").Value = -4752.5714

$ws.Range("H134").Value = 4007.3262
$ws.Range("I134").Value = 1966.1818
$ws.Range("K134").Value = 5898.5454
$ws.Range("M134").Value = -3363.5454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 312.41177
$ws.Range("I7").Value = 58.22222
$ws.Range("K7").Value = 58.22222
$ws.Range("M7").Value = 54.77778

$ws.Range("H17").Value = 20000
$ws.Range("I17").Value = 10000
$ws.Range("K17").Value = 10000
$ws.Range("M17").Value = -9826

$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()

$ws.Range("H25").Value = 7859
$ws.Range("I25").Value = 8073.75
$ws.Range("K25").Value = 8073.75
$ws.Range("M25").Value = -7899.75

$ws.Range("H31").Value = 5292.5713
$ws.Range("I31").Value = 2630.8125
$ws.Range("J31").Value = 6583.121
$ws.Range("K31").Value = 2630.8125
$ws.Range("L31").Value = 6583.121
$ws.Range("M31").Value = -2335.8125
$ws.Range("N31").Value = -7173.121

$ws.Range("H34").Value = 5292.5713
$ws.Range("I34").Value = 2630.8125
$ws.Range("J34").Value = 6583.121
$ws.Range("K34").Value = 2630.8125
$ws.Range("L34").Value = 6583.121
$ws.Range("M34").Value = -2428.8125
$ws.Range("N34").Value = -6987.121

$ws.Range("H88").Value = 16171
$ws.Range("J88").Value = 16171
$ws.Range("L88").Value = 16171
$ws.Range("N88").Value = -16983

$ws.Range("H91").Value = 16171
$ws.Range("J91").Value = 16171
$ws.Range("L91").Value = 16171
$ws.Range("N91").Value = -18979

$ws.Range("H94").Value = 1093.4546
$ws.Range("I94").Value = 737.4286
$ws.Range("J94").Value = 1259.6
$ws.Range("K94").Value = 737.4286
$ws.Range("L94").Value = 1259.6
$ws.Range("M94").Value = -286.4286
$ws.Range("N94").Value = -2161.6

$ws.Range("H134").Value = 1801
$ws.Range("I134").Value = 1801
$ws.Range("K134").Value = 5403
$ws.Range("M134").Value = -2868

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 10
$ws.Range("I6").Value = 10
$ws.Range("K6").Value = 30
$ws.Range("M6").Value = 83

$ws.Range("H12").Value = 20
$ws.Range("I12").Value = 15.166667
$ws.Range("J12").Value = 25.8
$ws.Range("K12").Value = 45.500001
$ws.Range("L12").Value = 77.4
$ws.Range("M12").Value = 127.499999
$ws.Range("N12").Value = -423.4

$ws.Range("H13").Value = 237.6
$ws.Range("I13").Value = 199
$ws.Range("K13").Value = 597
$ws.Range("M13").Value = -429

$ws.Range("H57").Value = 2478
$ws.Range("I57").Value = 1637.3334
$ws.Range("K57").Value = 4912.0002
$ws.Range("M57").Value = -4353.0002

$ws.Range("H129").Value = 1672.0588
$ws.Range("I129").Value = 892.36365
$ws.Range("J129").Value = 3101.5
$ws.Range("K129").Value = 2677.09095
$ws.Range("L129").Value = 9304.5
$ws.Range("M129").Value = 2322.90905
$ws.Range("N129").Value = -19304.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5490.7
$ws.Range("I122").Value = 2942.5833
$ws.Range("K122").Value = 8827.749899999999
$ws.Range("M122").Value = -6377.749899999999

$ws.Range("H132").Value = 1290.75
$ws.Range("I132").Value = 1317.1818
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 3951.5454
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -1421.5454
$ws.Range("N132").Value = -8060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3549.5
$ws.Range("I7").Value = 3132.8333
$ws.Range("K7").Value = 3132.8333
$ws.Range("M7").Value = -3020.8333

$ws.Range("H16").Value = 6774.125
$ws.Range("I16").Value = 9610.272
$ws.Range("K16").Value = 9610.272
$ws.Range("M16").Value = -9440.272

$ws.Range("H22").Value = 113509.78
$ws.Range("J22").Value = 3298
$ws.Range("L22").Value = 3298
$ws.Range("N22").Value = -3888

$ws.Range("H27").Value = 113509.78
$ws.Range("J27").Value = 3298
$ws.Range("L27").Value = 3298
$ws.Range("N27").Value = -3512

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").ClearContents()

$ws.Range("H61").Value = 2336.68
$ws.Range("I61").Value = 1618.3529
$ws.Range("K61").Value = 1618.3529
$ws.Range("M61").Value = -1416.3529

$ws.Range("H93").Value = 3184.75
$ws.Range("I93").Value = 3660.6667
$ws.Range("K93").Value = 3660.6667
$ws.Range("M93").Value = -2412.6667

$ws.Range("H100").Value = 7662.2915
$ws.Range("I100").Value = 2930.3845
$ws.Range("J100").Value = 13254.546
$ws.Range("K100").Value = 2930.3845
$ws.Range("L100").Value = 13254.546
$ws.Range("M100").Value = -2389.3845
$ws.Range("N100").Value = -14336.546

$ws.Range("H113").Value = 2336.68
$ws.Range("I113").Value = 1618.3529
$ws.Range("K113").Value = 1618.3529
$ws.Range("M113").Value = 551.6470999999999

$ws.Range("H122").Value = 3563.7693
$ws.Range("I122").Value = 3291.0605
$ws.Range("K122").Value = 9873.1815
$ws.Range("M122").Value = -7423.181500000001

$ws.Range("H126").Value = 3549.5
$ws.Range("I126").Value = 3132.8333
$ws.Range("K126").Value = 9398.499899999999
$ws.Range("M126").Value = -6928.499899999999

$ws.Range("H128").Value = 93209.664
$ws.Range("J128").Value = 93209.664
$ws.Range("L128").Value = 93209.664
$ws.Range("N128").Value = -103169.664

$ws.Range("H132").Value = 2948.76
$ws.Range("I132").Value = 2948.76
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8846.28
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -6316.280000000001
$ws.Range("N132").ClearContents()

$ws.Range("H136").Value = 4799.6587
$ws.Range("I136").Value = 4156.4863
$ws.Range("J136").Value = 10749
$ws.Range("K136").Value = 12469.4589
$ws.Range("L136").Value = 32247
$ws.Range("M136").Value = -9919.458899999998
$ws.Range("N136").Value = -37347

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

$ws.Range("H41").Value = 15000
$ws.Range("J41").Value = 15000
$ws.Range("L41").Value = 15000
$ws.Range("N41").Value = -15780

$ws.Range("H96").Value = 1999.6
$ws.Range("I96").Value = 2199.5
$ws.Range("K96").Value = 2199.5
$ws.Range("M96").Value = -826.5

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H122").Value = 2445.2727
$ws.Range("J122").Value = 2962.375
$ws.Range("L122").Value = 8887.125
$ws.Range("N122").Value = -13787.125

$ws.Range("H126").Value = 7128.2856
$ws.Range("I126").Value = 6699.8
$ws.Range("J126").Value = 8199.5
$ws.Range("K126").Value = 20099.4
$ws.Range("L126").Value = 24598.5
$ws.Range("M126").Value = -17629.4
$ws.Range("N126").Value = -29538.5

$ws.Range("H132").Value = 3865.5615
$ws.Range("I132").Value = 3199.75
$ws.Range("K132").Value = 9599.25
$ws.Range("M132").Value = -7069.25

$ws.Range("H136").Value = 8038712
$ws.Range("I136").Value = 9874998
$ws.Range("J136").Value = 4961.875
$ws.Range("K136").Value = 29624994
$ws.Range("L136").Value = 14885.625
$ws.Range("M136").Value = -29622444
$ws.Range("N136").Value = -19985.625

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
